$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1879.134
$ws.Range("I15").Value = 1879.134
$ws.Range("K15").Value = 5637.402
$ws.Range("M15").Value = -5468.402
$ws.Range("H18").Value = 749.75
$ws.Range("I18").Value = 749.75
$ws.Range("K18").Value = 749.75
$ws.Range("M18").Value = -465.75
$ws.Range("H19").Value = 1021.5
$ws.Range("I19").Value = 795.0833
$ws.Range("J19").Value = 1293.2
$ws.Range("K19").Value = 795.0833
$ws.Range("L19").Value = 1293.2
$ws.Range("M19").Value = -620.0833
$ws.Range("N19").Value = -1643.2
$ws.Range("H51").Value = 6123.3335
$ws.Range("I51").Value = 3116
$ws.Range("J51").Value = 8271.429
$ws.Range("K51").Value = 3116
$ws.Range("L51").Value = 8271.429
$ws.Range("M51").Value = -2632
$ws.Range("N51").Value = -9239.429
$ws.Range("H70").Value = 2775
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2775
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 8325
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -8865
$ws.Range("H73").Value = 2775
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2775
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 8325
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -10197
$ws.Range("H125").Value = 3097.1904
$ws.Range("I125").Value = 3266.8
$ws.Range("J125").Value = 3044.1875
$ws.Range("K125").Value = 29401.2
$ws.Range("L125").Value = 27397.6875
$ws.Range("M125").Value = -26941.2
$ws.Range("N125").Value = -32317.6875
$ws.Range("H129").Value = 1060.6389
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1060.6389
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3181.9167
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -13181.9167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2884
$ws.Range("I63").Value = 2127.1428
$ws.Range("J63").Value = 4650
$ws.Range("K63").Value = 2127.1428
$ws.Range("L63").Value = 4650
$ws.Range("M63").Value = -1441.1428
$ws.Range("N63").Value = -6022
$ws.Range("H66").Value = 2884
$ws.Range("I66").Value = 2127.1428
$ws.Range("J66").Value = 4650
$ws.Range("K66").Value = 10635.714
$ws.Range("L66").Value = 23250
$ws.Range("M66").Value = -7203.714
$ws.Range("N66").Value = -30114
$ws.Range("H105").Value = 28390
$ws.Range("J105").Value = 28390
$ws.Range("L105").Value = 28390
$ws.Range("N105").Value = -35378

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 26000
$ws.Range("J44").Value = 26000
$ws.Range("L44").Value = 26000
$ws.Range("N44").Value = -26994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2530.5757
$ws.Range("I31").Value = 1364.3636
$ws.Range("J31").Value = 3113.682
$ws.Range("K31").Value = 1364.3636
$ws.Range("L31").Value = 3113.682
$ws.Range("M31").Value = -1069.3636
$ws.Range("N31").Value = -3703.682
$ws.Range("H33").Value = 15540.667
$ws.Range("J33").Value = 30508.75
$ws.Range("L33").Value = 30508.75
$ws.Range("N33").Value = -31266.75
$ws.Range("H34").Value = 2530.5757
$ws.Range("I34").Value = 1364.3636
$ws.Range("J34").Value = 3113.682
$ws.Range("K34").Value = 1364.3636
$ws.Range("L34").Value = 3113.682
$ws.Range("M34").Value = -1162.3636
$ws.Range("N34").Value = -3517.682
$ws.Range("H122").Value = 1199
$ws.Range("I122").Value = 1154
$ws.Range("J122").Value = 1266.5
$ws.Range("K122").Value = 3462
$ws.Range("L122").Value = 3799.5
$ws.Range("M122").Value = -1012
$ws.Range("N122").Value = -8699.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1045497.75
$ws.Range("I113").Value = 2020766
$ws.Range("J113").Value = 567.4286
$ws.Range("K113").Value = 6062298
$ws.Range("L113").Value = 1702.2858
$ws.Range("M113").Value = -6060128
$ws.Range("N113").Value = -6042.2858
$ws.Range("H131").Value = 951.26
$ws.Range("J131").Value = 983.86316
$ws.Range("L131").Value = 2951.58948
$ws.Range("N131").Value = -13031.58948
$ws.Range("H134").Value = 6180.3335
$ws.Range("I134").Value = 3750
$ws.Range("J134").Value = 8610.666999999999
$ws.Range("K134").Value = 11250
$ws.Range("L134").Value = 25832.001
$ws.Range("M134").Value = -6180
$ws.Range("N134").Value = -35972.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7162.077
$ws.Range("I43").Value = 821.4
$ws.Range("J43").Value = 11125
$ws.Range("K43").Value = 821.4
$ws.Range("L43").Value = 11125
$ws.Range("M43").Value = -670.4
$ws.Range("N43").Value = -11427
$ws.Range("H49").Value = 18833.334
$ws.Range("J49").Value = 18833.334
$ws.Range("L49").Value = 18833.334
$ws.Range("N49").Value = -19201.334
$ws.Range("H107").Value = 838.17645
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 1094.3334
$ws.Range("K107").Value = 550
$ws.Range("L107").Value = 1094.3334
$ws.Range("M107").Value = 1370
$ws.Range("N107").Value = -4934.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 948.8421
$ws.Range("I46").Value = 1010.1
$ws.Range("J46").Value = 880.7778
$ws.Range("K46").Value = 1010.1
$ws.Range("L46").Value = 880.7778
$ws.Range("M46").Value = -822.1
$ws.Range("N46").Value = -1256.7778
$ws.Range("H132").Value = 2758.3257
$ws.Range("I132").Value = 2663.6191
$ws.Range("J132").Value = 2848.7273
$ws.Range("K132").Value = 7990.8573
$ws.Range("L132").Value = 8546.1819
$ws.Range("M132").Value = -5460.8573
$ws.Range("N132").Value = -13606.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -21144
$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = -16138
$ws.Range("H126").Value = 1025
$ws.Range("I126").Value = 971.4286
$ws.Range("K126").Value = 2914.2858
$ws.Range("M126").Value = -444.2857999999997
